# "Updated symbol list on Tue Dec 20 07:43:55 UTC 2022 with GitHub Actions"
#
# The sheet stores every cell (besides column A) as text, including the
# "Price" column (D), which holds numeric-looking strings such as
# "248.31". A plain Range.Value/.Formula assignment of a numeric-looking
# string is auto-coerced to a real number by the engine (like typing it
# into Excel), so those writes are prefixed with a leading apostrophe to
# force them to stay text, exactly like the source data. Plain text
# columns (B/C/E - coin name, link, summary) don't need the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Formula = "'248.16"

# Row 3 - OKB
$ws.Range("D3").Formula = "'21.67"

# Row 4 - HuobiToken
$ws.Range("D4").Formula = "'5.307"

# Row 5 - Cronos
$ws.Range("D5").Formula = "'0.05615"

# Row 6
$ws.Range("D6").Formula = "'3.418"

# Row 7
$ws.Range("D7").Formula = "'6.388"

# Row 8
$ws.Range("D8").Formula = "'0.8122"

# Row 9
$ws.Range("D9").Formula = "'0.9472"

# Row 10
$ws.Range("D10").Formula = "'0.1428"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").Formula = "'0.07510"

# Row 13 - BitrueCoin
$ws.Range("D13").Formula = "'0.03096"

# Row 14
$ws.Range("D14").Formula = "'0.09309"

# Row 15
$ws.Range("D15").Formula = "'3.586"

# Row 16
$ws.Range("D16").Formula = "'0.001588"

# Row 17
$ws.Range("D17").Formula = "'0.04713"

# Row 18 - One
$ws.Range("D18").Formula = "'0.0005782"
$ws.Range("E18").Formula = "17OneONEWorstin24h"

# Row 19
$ws.Range("D19").Formula = "'0.006311"

# Row 20
$ws.Range("D20").Formula = "'0.005047"

# Row 21
$ws.Range("D21").Formula = "'0.001034"

# Row 22
$ws.Range("D22").Formula = "'0.0001501"

# Row 23
$ws.Range("D23").Formula = "'3.778"

# Row 25
$ws.Range("D25").Formula = "'0.3302"

# Row 26
$ws.Range("D26").Formula = "'0.1308"

# Row 28
$ws.Range("D28").Formula = "'0.0003001"

# Row 40 - IDEX
$ws.Range("D40").Formula = "'0.03955"

# Row 41 - was BKEXToken, now KickToken (rows 41-43 rotate: BKEXToken,
# CEJI, KickToken -> KickToken, BKEXToken, CEJI)
$ws.Range("B41").Formula = "KickToken"
$ws.Range("C41").Formula = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Formula = "'0.006864"
$ws.Range("E41").Formula = "40KickTokenKICK"

# Row 42 - was CEJI, now BKEXToken
$ws.Range("B42").Formula = "BKEXToken"
$ws.Range("C42").Formula = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Formula = "'0.1061"
$ws.Range("E42").Formula = "41BKEXTokenBKK"

# Row 43 - was KickToken, now CEJI
$ws.Range("B43").Formula = "CEJI"
$ws.Range("C43").Formula = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Formula = "'0.003113"
$ws.Range("E43").Formula = "42CEJICEJI"

# Row 44 - LocalTraders
$ws.Range("D44").Formula = "'0.008799"

# Row 45 - CoinLion
$ws.Range("D45").Formula = "'0.00005603"

# Row 47 - ACDXExchange
$ws.Range("D47").Formula = "'0.0005502"
$ws.Range("E47").Formula = "46ACDXExchangeACXT"

# Row 48 - CoinbaseStockToken
$ws.Range("D48").Formula = "'0.7803"

# Row 49 - BOLO
$ws.Range("D49").Formula = "'0.1773"

# Row 50 - CryptobidCoin
$ws.Range("D50").Formula = "'0.00002101"

# Row 51 - SpecialPowerGold
$ws.Range("D51").Formula = "'0.01011"
